# Generate Report for Handoff
#
# For the six "Ready for handoff" rows (27ab8485, 7d469ff3, 8dc3976d,
# c9f9b0cf, e1516a6e, fdbcf554 -> rows 8,9,10,12,13,14) on both the
# zh-cn and de-de per-language sheets:
#   - set the Priority column (E) to "ht" (it was previously blank)
#   - bump the "Latest Handoff Datetime" column (H) forward a bit, as a
#     fresh handoff report was generated
# and mirror the later handoff timestamp on the Overview sheet's
# "Latest HO Xliff Generate Date" column (G) for the de-de-driven rows.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 12, 13, 14)

# Overview sheet: column G timestamp refresh (shared with de-de's column H).
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-18 22:21:48"
}

# zh-cn sheet: Priority -> "ht", handoff datetime refresh.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-18 22:21:42"
}

# de-de sheet: Priority -> "ht", handoff datetime refresh.
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-18 22:21:48"
}
